# Weighting & Scaling update & heatmap
# Applies the "Scaling" sheet changes:
#   - C2:C4 become MAX() formulas over ecological_params
#   - New "Optimal" / "Threshold" columns (E, F) with header + bordered blank cells
#   - Selection / dimension follow along automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scaling")

# --- C2:C4 -> replace literal 100 with MAX() formulas referencing ecological_params ---
$ws.Range("C2").Formula = "=MAX(ecological_params!B2:D2)"
$ws.Range("C3").Formula = "=MAX(ecological_params!B3:D3)"
$ws.Range("C4").Formula = "=MAX(ecological_params!B4:D4)"

# --- New header cells E1 ("Optimal") and F1 ("Threshold"), formatted like the other headers ---
$ws.Range("C1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "Optimal"
$ws.Range("F1").Value = "Threshold"

# --- New blank data cells E2:F4, bordered (placeholders for the heatmap) ---
$dataCells = $ws.Range("E2:F4")
$dataCells.Borders.Color = 0
$dataCells.Borders.LineStyle = 1

# --- Update the sheet's remembered selection ---
$ws.Activate()
$ws.Range("C10").Select()

# Restore the original active sheet so the workbook's active tab is unchanged
$wb.Worksheets.Item("ecological_params").Activate()
